$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7219
$ws1.Range("F3").Value = 60
$ws1.Range("F5").Value = 157
$ws1.Range("F6").Value = 1103
$ws1.Range("F7").Value = 181
$ws1.Range("F8").Value = 8
$ws1.Range("F9").Value = 76
$ws1.Range("F10").Value = 16

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7219
$ws4.Range("F3").Value = 60
$ws4.Range("F5").Value = 157
$ws4.Range("F6").Value = 1103
$ws4.Range("F7").Value = 181
$ws4.Range("F9").Value = 8
$ws4.Range("F10").Value = 76
$ws4.Range("F11").Value = 16

$wb.Save()
